# Updates cryptos list prices / volumes / reorders a few coin rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "59.162.10"
$ws.Range("E2").Value = "  -6.10%  "

# Row 3
$ws.Range("D3").Value = "2.468.53"
$ws.Range("E3").Value = "  -8.18%  "

# Row 4
$ws.Range("E4").Value = "  -0.19%  "

# Row 5
$ws.Range("D5").Value = "'541.66"
$ws.Range("E5").Value = "  -2.67%  "

# Row 6
$ws.Range("D6").Value = "'147.40"
$ws.Range("E6").Value = "  -7.27%  "

# Row 7
$ws.Range("E7").Value = "  -0.20%  "

# Row 8
$ws.Range("D8").Value = "'0.571"
$ws.Range("E8").Value = "  -3.72%  "

# Row 9
$ws.Range("D9").Value = "2.465.32"
$ws.Range("E9").Value = "  -8.56%  "

# Row 10
$ws.Range("E10").Value = "  -6.51%  "

# Row 11
$ws.Range("D11").Value = "'0.159"
$ws.Range("E11").Value = "  -2.15%  "

# Row 12
$ws.Range("D12").Value = "'5.34"
$ws.Range("E12").Value = "  -0.72%  "

# Row 13
$ws.Range("E13").Value = "  -5.00%  "

# Row 14
$ws.Range("D14").Value = "2.895.89"
$ws.Range("E14").Value = "  -8.54%  "

# Row 15
$ws.Range("D15").Value = "'24.07"
$ws.Range("E15").Value = "  -9.91%  "

# Row 16
$ws.Range("D16").Value = "59.042.30"
$ws.Range("E16").Value = "  -6.18%  "

# Row 17
$ws.Range("E17").Value = "  -6.65%  "

# Row 18
$ws.Range("D18").Value = "2.512.90"
$ws.Range("E18").Value = "  -6.58%  "

# Row 19
$ws.Range("E19").Value = "  -6.78%  "

# Row 21
$ws.Range("D21").Value = "'324.96"
$ws.Range("E21").Value = "  -6.14%  "

# Row 22
$ws.Range("E22").Value = "  -3.30%  "

# Row 23
$ws.Range("D23").Value = "'5.72"
$ws.Range("E23").Value = "  -9.62%  "

# Row 24
$ws.Range("E24").Value = "  -10.16%  "

# Row 25
$ws.Range("D25").Value = "'60.69"
$ws.Range("E25").Value = "  -4.60%  "

# Row 26
$ws.Range("D26").Value = "'0.160"
$ws.Range("E26").Value = "  -5.03%  "

# Row 27
$ws.Range("D27").Value = "'0.979"
$ws.Range("E27").Value = "  -1.93%  "

# Row 28
$ws.Range("D28").Value = "'7.73"
$ws.Range("E28").Value = "  -6.53%  "

# Row 29
$ws.Range("D29").Value = "'1.28"
$ws.Range("E29").Value = "  -11.79%  "

# Row 30
$ws.Range("E30").Value = "  -6.38%  "

# Row 31
$ws.Range("E31").Value = "  -10.38%  "

# Row 32
$ws.Range("D32").Value = "'6.66"
$ws.Range("E32").Value = "  -9.30%  "

# Row 33
$ws.Range("D33").Value = "'0.997"
$ws.Range("E33").Value = "  -0.11%  "

# Row 34
$ws.Range("D34").Value = "'157.26"
$ws.Range("E34").Value = "  -4.38%  "

# Row 35
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'1.37"
$ws.Range("E35").Value = "  -9.47%  "

# Row 36
$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").Value = "'18.42"
$ws.Range("E36").Value = "  -5.91%  "

# Row 37
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").Value = "'4.48"
$ws.Range("E37").Value = "  -10.12%  "

# Row 38
$ws.Range("D38").Value = "'1.72"
$ws.Range("E38").Value = "  -4.20%  "

# Row 39
$ws.Range("D39").Value = "'320.15"
$ws.Range("E39").Value = "  -11.48%  "

# Row 40
$ws.Range("D40").Value = "'5.90"
$ws.Range("E40").Value = "  -9.26%  "

# Row 41
$ws.Range("D41").Value = "'36.30"
$ws.Range("E41").Value = "  -5.70%  "

# Row 42
$ws.Range("E42").Value = "  -13.72%  "

# Row 43
$ws.Range("D43").Value = "'3.70"
$ws.Range("E43").Value = "  -8.15%  "

# Row 44
$ws.Range("D44").Value = "'0.997"
$ws.Range("E44").Value = "  -0.11%  "

# Row 45
$ws.Range("D45").Value = "'10.73"
$ws.Range("E45").Value = "  -2.71%  "

# Row 46
$ws.Range("D46").Value = "'0.0943"
$ws.Range("E46").Value = "  -3.23%  "

# Row 47
$ws.Range("D47").Value = "'0.582"
$ws.Range("E47").Value = "  -6.31%  "

# Row 48
$ws.Range("E48").Value = "  -6.77%  "

# Row 49
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "'0.0230"
$ws.Range("E49").Value = "  -6.07%  "

# Row 50
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "'18.99"
$ws.Range("E50").Value = "  -10.43%  "

# Row 51
$ws.Range("D51").Value = "'121.90"
$ws.Range("E51").Value = "  -5.60%  "
